$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record was added to the spreadsheet (ID 8: Ryan Rain).
# This mirrors the row that should trigger the stored procedure call.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Ryan"
$ws.Range("C9").Value = "Rain"

$ws.Range("F9").Select()
